$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline")

# --- Update row 5 ("bert + lstm") with refreshed training metrics ---
$ws.Range("C5").Value = 1.4444999999999999
$ws.Range("D5").Value = 0.67720000000000002
$ws.Range("E5").Value = 0.67520000000000002
$ws.Range("F5").Value = 0.67720000000000002
$ws.Range("G5").Value = 0.66369999999999996
$ws.Range("H5").Value = 0.374
$ws.Range("I5").Value = 0.38100000000000001
$ws.Range("J5").Value = 0.57569999999999999

# --- New row 13: "deberta + bilstm" results (written first so its shared
#     string lands at the same index the reference workbook uses) ---
$ws.Range("B13").Value = "deberta + bilstm"
$ws.Range("C13:J13").NumberFormat = "0.0000"
$ws.Range("C13").Value = 2.4900000000000002
$ws.Range("D13").Value = 0.47
$ws.Range("E13").Value = 0.2258
$ws.Range("F13").Value = 0.47
$ws.Range("G13").Value = 0.30330000000000001
$ws.Range("H13").Value = 0.0145
$ws.Range("I13").Value = 0.0213
$ws.Range("J13").Value = 0.0227

# --- New column K: per-run training time, plus header ---
$ws.Range("K5").Value = "9 hrs 44 mins 30 secs"
$ws.Range("K1").Value = "total time"
$ws.Columns.Item(11).ColumnWidth = 19.43

$ws.Range("E6").Select()
